# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-28, replacing the previous
# Strike# derived values with the recalculated K values.
$kValues = @{
    2  = 1
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 4
    10 = 0
    11 = 0
    12 = 2
    13 = 0
    14 = 0
    15 = 1
    16 = 1
    17 = 0
    18 = 4
    19 = 2
    20 = 1
    21 = 0
    22 = 1
    23 = 0
    24 = 1
    25 = 0
    26 = 1
    27 = 7
    28 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
